$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.003.88"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").Value = "2.637.83"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.54"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.34"
$ws.Range("E6").Value = "  +3.59%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("E9").Value = "  +5.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.400"
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  +1.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.08"
$ws.Range("E13").Value = "  +5.43%  "
$ws.Range("E14").Value = "  +20.51%  "
$ws.Range("D15").Value = "3.110.71"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "64.883.44"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.602.58"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  +2.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.79"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "352.08"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.32"
$ws.Range("E21").Value = "  +6.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.01"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.94"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.54"
$ws.Range("E25").Value = "  +4.62%  "
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.12"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "0.0₃0939"
$ws.Range("E30").Value = "  +9.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.08"
$ws.Range("E31").Value = "  +2.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "505.12"
$ws.Range("E32").Value = "  -8.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.60"
$ws.Range("E34").Value = "  +6.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.31"
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.01"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.19"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("E39").Value = "  +5.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.35"
$ws.Range("E42").Value = "  +6.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "165.18"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0612"
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("E47").Value = "  +3.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.645"
$ws.Range("E48").Value = "  +2.73%  "
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0980"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.37"
$ws.Range("E51").Value = "  +0.89%  "
